# Generate Report for Handoff
# - Refresh the "Latest Handoff / HO Xliff Generate Date" timestamps for the
#   entries that were (re-)handed off in this run.
# - Mark those same entries' Priority as "ht" on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = 7,8,9,10,12,14

$overviewStamp = "2016-08-28 08:22:04"
$zhcnStamp     = "2016-08-28 08:21:57"
$dedeStamp     = "2016-08-28 08:22:04"

foreach ($r in $rows) {
    $overview.Range("G$r").Value = $overviewStamp

    $zhcn.Range("H$r").Value = $zhcnStamp
    $zhcn.Range("E$r").Value = "ht"

    $dede.Range("H$r").Value = $dedeStamp
    $dede.Range("E$r").Value = "ht"
}
